$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D timestamps: values cascade down (each block takes on the
# value previously held by the block above it), and a brand new timestamp
# is written to the first block - this mirrors an "Actualizar" (refresh)
# pass that re-stamped each existing entry with the timestamp of the row
# processed after it, and stamped the newest entries with the current run
# time.

$ws.Range("D2:D15").Value = 44244.52730866284
$ws.Range("D16:D29").Value = 44244.50605336806
$ws.Range("D30:D43").Value = 44244.48480421296
